# Refresh cryptocurrency price/volume data (and two re-ranked rows) to
# match the latest coinranking.com snapshot.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "64.869.33"
$ws.Range("E2").Value = "  -2.14%  "

$ws.Range("D3").Value = "3.234.85"
$ws.Range("E3").Value = "  -1.33%  "

$ws.Range("E4").Value = "  -0.04%  "

$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "578.25"
$ws.Range("E5").Value = "  -0.30%  "

$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "173.07"
$ws.Range("E6").Value = "  -3.38%  "

$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = "0.632"
$ws.Range("E7").Value = "  +0.38%  "

$ws.Range("E8").Value = "  -0.06%  "

$ws.Range("D9").Value = "3.236.21"
$ws.Range("E9").Value = "  -1.27%  "

$ws.Range("E10").Value = "  -2.64%  "

$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = "6.79"
$ws.Range("E11").Value = "  +0.67%  "

$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = "0.390"
$ws.Range("E12").Value = "  -2.91%  "

$ws.Range("D13").Value = "3.790.65"
$ws.Range("E13").Value = "  -1.55%  "

$ws.Range("E14").Value = "  -3.16%  "

$ws.Range("D15").Value = "64.914.41"
$ws.Range("E15").Value = "  -2.06%  "

$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = "25.81"
$ws.Range("E16").Value = "  -1.96%  "

$ws.Range("E17").Value = "  -2.64%  "

$ws.Range("D18").Value = "3.199.51"
$ws.Range("E18").Value = "  -1.98%  "

$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = "418.41"
$ws.Range("E19").Value = "  -3.80%  "

$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = "5.40"
$ws.Range("E20").Value = "  -2.14%  "

$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = "12.85"
$ws.Range("E21").Value = "  -2.34%  "

$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "7.22"
$ws.Range("E22").Value = "  -2.29%  "

$ws.Range("E23").Value = "  -0.07%  "

$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = "70.98"
$ws.Range("E24").Value = "  -1.08%  "

$ws.Range("E25").Value = "  -0.36%  "

$ws.Range("E26").Value = "  +4.42%  "

$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = "0.496"
$ws.Range("E27").Value = "  -1.75%  "

$ws.Range("E28").Value = "  -1.27%  "

$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = "9.12"
$ws.Range("E29").Value = "  +3.13%  "

$ws.Range("E30").Value = "  -0.09%  "

$ws.Range("E31").Value = "  -3.63%  "

$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = "21.84"
$ws.Range("E32").Value = "  -1.89%  "

$ws.Range("E33").Value = "  +0.03%  "

$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = "5.02"
$ws.Range("E34").Value = "  -3.32%  "

$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = "6.44"
$ws.Range("E35").Value = "  -2.07%  "

$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = "1.16"
$ws.Range("E36").Value = "  -1.29%  "

$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = "157.55"
$ws.Range("E37").Value = "  +0.05%  "

$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = "1.40"
$ws.Range("E38").Value = "  -1.73%  "

$ws.Range("D39").Value = "2.824.93"
$ws.Range("E39").Value = "  +1.83%  "

$ws.Range("E40").Value = "  -2.60%  "

$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = "25.45"
$ws.Range("E41").Value = "  -3.99%  "

$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = "4.24"
$ws.Range("E42").Value = "  -1.40%  "

$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = "39.62"
$ws.Range("E43").Value = "  -1.59%  "

$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = "0.726"
$ws.Range("E44").Value = "  -6.14%  "

$ws.Range("B45").Value = "RenderToken"
$ws.Range("C45").Value = "https://coinranking.com/coin/7C4Mh4xy1yDel+rendertoken-rndr"
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = "5.76"
$ws.Range("E45").Value = "  -4.56%  "

$ws.Range("B46").Value = "Hedera"
$ws.Range("C46").Value = "https://coinranking.com/coin/jad286TjB+hedera-hbar"
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = "0.0631"
$ws.Range("E46").Value = "  -4.02%  "

$ws.Range("B47").Value = "Bittensor"
$ws.Range("C47").Value = "https://coinranking.com/coin/pgv7xSFi6+bittensor-tao"
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = "303.36"
$ws.Range("E47").Value = "  -5.96%  "

$ws.Range("B48").Value = "dogwifhat"
$ws.Range("C48").Value = "https://coinranking.com/coin/sZUrmToWF+dogwifhat-wif"
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = "2.17"
$ws.Range("E48").Value = "  -4.92%  "

$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = "22.23"
$ws.Range("E49").Value = "  -4.07%  "

$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = "0.0265"
$ws.Range("E50").Value = "  -0.76%  "

$ws.Range("E51").Value = "  -1.04%  "
